# edit.ps1 - Applies the "Updated cryptos list on Sat Jan 27 18:13:55 UTC 2024
# with GitHub Actions" commit to the cryptos worksheet: refreshed prices and
# 1h volume deltas, plus a few coin rows that swapped rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry holds only the columns (B=Coin, C=Link, D=Price, E=Volume(1h))
# that actually changed for that row. "Text" marks Price values that must be
# forced to Text format so they don't get auto-converted to numbers (the
# sheet stores prices as literal strings, e.g. "41.758.00").
$updates = @(
    @{ Row=2; D='41.758.00'; E='  -0.60%  ' },
    @{ Row=3; D='2.267.04'; E='  -0.10%  ' },
    @{ Row=4; E='  +0.12%  ' },
    @{ Row=5; D='304.81'; DText=$true; E='  +0.81%  ' },
    @{ Row=6; D='92.78'; DText=$true; E='  -0.70%  ' },
    @{ Row=7; E='  -0.66%  ' },
    @{ Row=8; E='  +0.08%  ' },
    @{ Row=9; E='  -0.75%  ' },
    @{ Row=10; D='32.40'; DText=$true; E='  -1.48%  ' },
    @{ Row=11; D='0.0797'; DText=$true; E='  -0.69%  ' },
    @{ Row=12; D='0.112'; DText=$true; E='  -2.16%  ' },
    @{ Row=13; D='6.66'; DText=$true; E='  -0.77%  ' },
    @{ Row=14; D='2.617.87'; E='  -0.15%  ' },
    @{ Row=15; E='  +0.29%  ' },
    @{ Row=16; D='2.268.34'; E='  -0.24%  ' },
    @{ Row=17; D='0.782'; DText=$true; E='  +3.27%  ' },
    @{ Row=18; D='41.707.18'; E='  -0.44%  ' },
    @{ Row=19; D='12.75'; DText=$true; E='  +3.51%  ' },
    @{ Row=20; E='  -0.24%  ' },
    @{ Row=21; D='5.94'; DText=$true; E='  -0.17%  ' },
    @{ Row=22; D='67.57'; DText=$true; E='  +0.28%  ' },
    @{ Row=23; D='244.00'; DText=$true; E='  +0.48%  ' },
    @{ Row=24; D='2.60'; DText=$true; E='  +0.16%  ' },
    @{ Row=25; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='1.93'; DText=$true; E='  +1.16%  ' },
    @{ Row=26; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.00'; DText=$true; E='  +0.06%  ' },
    @{ Row=27; D='23.95'; DText=$true; E='  -0.05%  ' },
    @{ Row=28; D='9.58'; DText=$true; E='  -1.47%  ' },
    @{ Row=29; E='  -5.47%  ' },
    @{ Row=30; D='34.95'; DText=$true; E='  +2.34%  ' },
    @{ Row=31; D='159.18'; DText=$true; E='  +0.70%  ' },
    @{ Row=32; D='5.30'; DText=$true; E='  +2.16%  ' },
    @{ Row=33; E='  +0.08%  ' },
    @{ Row=34; D='0.0742'; DText=$true; E='  -0.41%  ' },
    @{ Row=35; D='3.03'; DText=$true; E='  -1.75%  ' },
    @{ Row=36; B='WEMIXToken'; C='https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D='2.37'; DText=$true; E='  -1.10%  ' },
    @{ Row=37; B='Celestia'; C='https://coinranking.com/coin/YQcD0lBl7+celestia-tia'; D='16.79'; DText=$true; E='  +1.11%  ' },
    @{ Row=38; E='  +0.68%  ' },
    @{ Row=39; E='  -0.24%  ' },
    @{ Row=40; E='  -0.19%  ' },
    @{ Row=41; D='3.92'; DText=$true; E='  -1.37%  ' },
    @{ Row=42; D='19.97'; DText=$true; E='  -1.40%  ' },
    @{ Row=43; D='2.007.98'; E='  -1.97%  ' },
    @{ Row=44; E='  +13.98%  ' },
    @{ Row=45; E='  +0.79%  ' },
    @{ Row=46; E='  +2.96%  ' },
    @{ Row=47; E='  -0.49%  ' },
    @{ Row=48; D='53.31'; DText=$true; E='  +2.90%  ' },
    @{ Row=49; D='73.12'; DText=$true; E='  +3.32%  ' },
    @{ Row=50; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.15'; DText=$true; E='  +0.56%  ' },
    @{ Row=51; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='1.50'; DText=$true; E='  -1.18%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey('B')) {
        $ws.Cells.Item($row, 2).Value = $u.B
    }
    if ($u.ContainsKey('C')) {
        $ws.Cells.Item($row, 3).Value = $u.C
    }
    if ($u.ContainsKey('D')) {
        $dCell = $ws.Cells.Item($row, 4)
        if ($u.ContainsKey('DText')) {
            # Value would otherwise be auto-recognised as a plain number
            # (e.g. "304.81"); force Text so it round-trips as a string.
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    if ($u.ContainsKey('E')) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
